$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Count -eq 2 -and $parts[0].Trim() -eq "dnasr281@gmail.com") {
            $cell.Value = $parts[1].Trim() + ", " + $parts[0].Trim()
        }
    }
}
